$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.268.26"
$ws.Range("E2").Value = "'  -0.22%  "
$ws.Range("D3").Value = "'1.904.63"
$ws.Range("E3").Value = "'  +0.18%  "
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("D5").Value = "'0.728"
$ws.Range("E5").Value = "'  +9.72%  "
$ws.Range("D6").Value = "'255.72"
$ws.Range("E6").Value = "'  +4.05%  "
$ws.Range("E7").Value = "'  -0.10%  "
$ws.Range("D8").Value = "'40.65"
$ws.Range("E8").Value = "'  -1.51%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("E9").Value = "'  +7.51%  "
$ws.Range("D10").Value = "'52.83"
$ws.Range("E10").Value = "'  -0.25%  "
$ws.Range("D11").Value = "'0.0761"
$ws.Range("E11").Value = "'  +5.89%  "
$ws.Range("D12").Value = "'0.0989"
$ws.Range("E12").Value = "'  -0.61%  "
$ws.Range("D13").Value = "'2.183.26"
$ws.Range("E13").Value = "'  +0.27%  "
$ws.Range("D14").Value = "'12.97"
$ws.Range("E14").Value = "'  +7.14%  "
$ws.Range("D15").Value = "'0.726"
$ws.Range("E15").Value = "'  +4.39%  "
$ws.Range("D16").Value = "'4.96"
$ws.Range("E16").Value = "'  +2.60%  "
$ws.Range("D17").Value = "'1.901.50"
$ws.Range("E17").Value = "'  -0.01%  "
$ws.Range("D18").Value = "'35.270.59"
$ws.Range("E18").Value = "'  -0.12%  "
$ws.Range("D19").Value = "'74.72"
$ws.Range("E19").Value = "'  +3.45%  "
$ws.Range("E20").Value = "'  +3.90%  "
$ws.Range("D21").Value = "'243.67"
$ws.Range("E21").Value = "'  +1.26%  "
$ws.Range("D22").Value = "'13.03"
$ws.Range("E22").Value = "'  +4.86%  "
$ws.Range("D23").Value = "'5.11"
$ws.Range("E23").Value = "'  +5.51%  "
$ws.Range("E24").Value = "'  -0.06%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "'  +7.45%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "'  +4.30%  "
$ws.Range("D27").Value = "'166.19"
$ws.Range("E27").Value = "'  -2.30%  "
$ws.Range("D28").Value = "'8.69"
$ws.Range("E28").Value = "'  +3.56%  "
$ws.Range("D29").Value = "'18.75"
$ws.Range("E29").Value = "'  +2.30%  "
$ws.Range("E30").Value = "'  +4.64%  "
$ws.Range("D31").Value = "'4.128.97"
$ws.Range("E31").Value = "'  +19.46%  "
$ws.Range("D32").Value = "'4.38"
$ws.Range("E32").Value = "'  +6.19%  "
$ws.Range("D33").Value = "'1.99"
$ws.Range("E33").Value = "'  +14.44%  "
$ws.Range("D34").Value = "'1.65"
$ws.Range("E34").Value = "'  +23.36%  "
$ws.Range("E35").Value = "'  +4.45%  "
$ws.Range("D36").Value = "'4.25"
$ws.Range("E36").Value = "'  +4.41%  "
$ws.Range("E37").Value = "'  -0.16%  "
$ws.Range("D38").Value = "'0.909"
$ws.Range("D39").Value = "'2.03"
$ws.Range("E39").Value = "'  +0.42%  "
$ws.Range("D40").Value = "'0.0218"
$ws.Range("E40").Value = "'  +5.17%  "
$ws.Range("D41").Value = "'17.08"
$ws.Range("E41").Value = "'  +6.34%  "
$ws.Range("D42").Value = "'96.40"
$ws.Range("E42").Value = "'  +7.71%  "
$ws.Range("E43").Value = "'  +1.94%  "
$ws.Range("D44").Value = "'0.0649"
$ws.Range("E44").Value = "'  +3.65%  "
$ws.Range("D45").Value = "'1.336.29"
$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("D46").Value = "'2.44"
$ws.Range("E46").Value = "'  +3.05%  "
$ws.Range("E47").Value = "'  +1.22%  "
$ws.Range("D48").Value = "'6.70"
$ws.Range("E48").Value = "'  +3.30%  "
$ws.Range("D49").Value = "'2.75"
$ws.Range("D50").Value = "'45.12"
$ws.Range("E50").Value = "'  -5.97%  "
$ws.Range("D51").Value = "'0.0755"
$ws.Range("E51").Value = "'  +7.12%  "
